$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 215 (shifts existing rows 215-236 down to 216-237)
$ws.Rows.Item(215).Insert()

# Populate the newly inserted row 215 with the new price-report record
$ws.Cells.Item(215, 1).Value = 3
$ws.Cells.Item(215, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(215, 3).Value = "Coquimbo"
$ws.Cells.Item(215, 4).Value = 44918
$ws.Cells.Item(215, 5).Value = 5
$ws.Cells.Item(215, 6).Value = 100112010
$ws.Cells.Item(215, 7).Value = "Achicoria"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 50
$ws.Cells.Item(215, 11).Value = 7000
$ws.Cells.Item(215, 12).Value = 7000
$ws.Cells.Item(215, 13).Value = 7000
$ws.Cells.Item(215, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(215, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(215, 16).Value = 438
$ws.Cells.Item(215, 17).Value = 16
$ws.Cells.Item(215, 18).Value = "Hortaliza"
